$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New content, written in the same order the original author must have
# entered it (our engine assigns shared-string indices in write order, and
# the target sharedStrings.xml shows indices 14-21 appearing in this order:
# puipui, G7-url, ASLONG..., G16-url, G6-url, G8-url, skateborad wheel, 套装).

# Row 7: "puipui" description + taobao link (G7)
$ws.Range("C7").Value = "puipui"
$g7Full = "https://item.taobao.com/item.htm?spm=a230r.1.14.75.63f46c1cSMRgZK&id=670320551279&ns=1&abbucket=13#detail"
$g7Addr = "https://item.taobao.com/item.htm?spm=a230r.1.14.75.63f46c1cSMRgZK&id=670320551279&ns=1&abbucket=13"
$ws.Hyperlinks.Add($ws.Range("G7"), $g7Addr, "detail", "", $g7Full)

# Row 16: motor listing text + plain (non-hyperlinked) tmall link text
$ws.Range("C16").Value = "ASLONG PG16-050行星减速电机 微型直流马达 智能锁 电动牙刷12V"
$ws.Range("G16").Value = "https://detail.tmall.com/item.htm?spm=a220o.1000855.w4004-15686600035.7.1ae2f64902fGTW&id=665631367082&skuId=4789727204925"

# Row 6: puipui filling link (G6)
$g6Full = "https://item.taobao.com/item.htm?spm=a230r.1.14.18.4367136ehOktAT&id=19027438685&ns=1&abbucket=13#detail"
$g6Addr = "https://item.taobao.com/item.htm?spm=a230r.1.14.18.4367136ehOktAT&id=19027438685&ns=1&abbucket=13"
$ws.Hyperlinks.Add($ws.Range("G6"), $g6Addr, "detail", "", $g6Full)

# Row 8: skateboard wheel link (G8)
$g8Full = "https://item.taobao.com/item.htm?id=566418069461&ali_refid=a3_430620_1006:1151032306:N:H5waukvccbOG93BXXZykyA%3D%3D:3dbe3d683fee580cbab33fdf9ef00abb&ali_trackid=1_3dbe3d683fee580cbab33fdf9ef00abb&spm=a230r.1.14.6#detail"
$g8Addr = "https://item.taobao.com/item.htm?id=566418069461&ali_refid=a3_430620_1006:1151032306:N:H5waukvccbOG93BXXZykyA%3D%3D:3dbe3d683fee580cbab33fdf9ef00abb&ali_trackid=1_3dbe3d683fee580cbab33fdf9ef00abb&spm=a230r.1.14.6"
$ws.Hyperlinks.Add($ws.Range("G8"), $g8Addr, "detail", "", $g8Full)

# Row 8: product name + colour/spec text
$ws.Range("C8").Value = "skateborad wheel"
$ws.Range("D8").Value = "6647磨砂粉色套装"

# Row 15: a lingering hyperlink-styled (but empty / unlinked) cell
$ws.Range("G15").Style = "Hyperlink"

# --- Column widths the user manually resized ---
$ws.Columns("C").ColumnWidth = 66.83333333333333
$ws.Columns("D").ColumnWidth = 25
$ws.Columns("F").ColumnWidth = 17.5
$ws.Columns("G").ColumnWidth = 23.333333333333332

# --- Selection moved to C5 ---
$ws.Range("C5").Select()
